# Applies the attendance_reports sync changes described in the commit:
# - Reorder recorder names in G2 (insert "System" before "Veronia.rafat...")
# - Reorder recorder names in G9 (Shimaa.ashraf now first)
# - Fix Average Attendance % in L10 (10.2% -> 18.2%)
# - Fix Average Attendance % in S15 (10.2% -> 18.2%)
# - Add "maryam.ashraf@med.asu.edu.eg" to the recorder list in G28
# - Update attendance count in H28 (6/251 -> 66/251)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G2: recorder list reorder ---
$ws.Range("G2").Value = "Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

# --- G9: recorder list reorder ---
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# --- L10: Average Attendance % value update ---
# The text "18.2%" looks like a percentage, so a plain .Value assignment would be
# auto-converted into a numeric percent value and pick up a new number format /
# style. Force it to stay plain text by briefly using a text number format, then
# restore the original look (fill/alignment/format) by copying it from a
# neighboring cell that already carries the unchanged style (s="5").
$cellL10 = $ws.Range("L10")
$cellL10.NumberFormat = "@"
$cellL10.Value = "18.2%"
$ws.Range("K10").Copy() | Out-Null
$cellL10.PasteSpecial(-4122) | Out-Null

# --- S15: Average Attendance % value update (same trick as L10) ---
$cellS15 = $ws.Range("S15")
$cellS15.NumberFormat = "@"
$cellS15.Value = "18.2%"
$ws.Range("R15").Copy() | Out-Null
$cellS15.PasteSpecial(-4122) | Out-Null

# --- G28: append additional recorder ---
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

# --- H28: updated attendance count ---
$ws.Range("H28").Value = "66/251"

$excel.CutCopyMode = 0
